$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.052.94"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "3.081.51"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'234.32"
$ws.Range("D6").Value = "'608.87"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("D7").Value = "'1.10"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("E8").Value = "  -4.82%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'0.828"
$ws.Range("E10").Value = "  +13.35%  "
$ws.Range("D11").Value = "3.081.07"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  -2.84%  "
$ws.Range("D13").Value = "93.943.66"
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D14").Value = "'0.0000241"
$ws.Range("E14").Value = "  -5.40%  "
$ws.Range("D15").Value = "'34.06"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "'5.27"
$ws.Range("E16").Value = "  -4.07%  "
$ws.Range("D17").Value = "3.663.85"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "3.117.91"
$ws.Range("E18").Value = "  +1.52%  "
$ws.Range("D19").Value = "'3.66"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Value = "'14.60"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "'5.77"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "'443.09"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").Value = "'8.84"
$ws.Range("E23").Value = "  -5.75%  "
$ws.Range("D24").Value = "'0.0000192"
$ws.Range("E24").Value = "  -5.70%  "
$ws.Range("D25").Value = "'8.29"
$ws.Range("E25").Value = "  +5.07%  "
$ws.Range("D26").Value = "'5.54"
$ws.Range("E26").Value = "  -4.06%  "
$ws.Range("D27").Value = "'84.94"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("D28").Value = "'11.90"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").Value = "3.252.64"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'0.251"
$ws.Range("E31").Value = "  +8.65%  "
$ws.Range("D32").Value = "'0.179"
$ws.Range("E32").Value = "  +5.87%  "
$ws.Range("D33").Value = "'0.124"
$ws.Range("E33").Value = "  -9.90%  "
$ws.Range("D34").Value = "'9.23"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").Value = "'7.75"
$ws.Range("E36").Value = "  -4.18%  "
$ws.Range("D37").Value = "'0.158"
$ws.Range("E37").Value = "  -3.58%  "
$ws.Range("D38").Value = "'25.65"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("D40").Value = "'0.446"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").Value = "'23.95"
$ws.Range("D42").Value = "'1.27"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").Value = "'467.42"
$ws.Range("E43").Value = "  -3.39%  "
$ws.Range("D44").Value = "'3.69"
$ws.Range("E44").Value = "  -13.61%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'3.14"
$ws.Range("E46").Value = "  -10.55%  "
$ws.Range("D47").Value = "'161.50"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").Value = "'1.85"
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("D49").Value = "'0.677"
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("D50").Value = "'43.70"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").Value = "'0.998"
$ws.Range("E51").Value = "  -0.01%  "
